$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1-5 V")

# --- New "Div coeff table" in columns E:F (rows 10-18) ---

# Header row: "V" label over both value/result columns
$ws.Range("E10").Value = "V"
$ws.Range("F10").Value = "V"

# Multiplier values in E11:E15 and their corresponding results in F11:F15
# (plus the trailing spacer/checking rows), all formatted with the new
# "0.000" number format.
$ws.Range("E11:F18").NumberFormat = "0.000"

$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2.58
$ws.Range("E14").Value = 4
$ws.Range("E15").Value = 5

$ws.Range("F11").Formula = "=E11*B$11"
$ws.Range("F12").Formula = "=E12*B$11"
$ws.Range("F13").Formula = "=E13*B$11"
$ws.Range("F14").Formula = "=E14*B$11"
$ws.Range("F15").Formula = "=E15*B$11"

# Checking row: recompute the divider ratio directly (bug fix - was showing 0.00 MPa)
$ws.Range("E18").Formula = "=B1/(B1+B2)"

# Make "1-5 V" the active sheet/tab and select E13, matching the saved view state
$ws.Activate()
$ws.Range("E13").Select()
